$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / info block updates ---
# Supplier name
$ws.Range("B8").Value = "Cameron Nieve"
# IAR No. (blanked out, keep a single space)
$ws.Range("E8").Value = " "
# PO No./Date
$ws.Range("B9").Value = " 2022-099283"
# Date
$ws.Range("E9").Value = "1/10/2023"
# Invoice No.
$ws.Range("B10").Value = " 12345678"
# Respo Center Code
$ws.Range("B11").Value = "AMO 310200100000"
# Date next to Respo Center Code (was duplicating the "Date" value, now its own date)
$ws.Range("E11").Value = "12/21/2022"

# --- Line item rows ---
$ws.Range("B15").Value = "as scdsfd safsdfdf dsaf"
$ws.Range("D15").Value = "(pcs) 1,000.00"
$ws.Range("E14").Value = 5
$ws.Range("E15").Value = 5

# --- Inspection / acceptance dates ---
$ws.Range("A18").Value = "Date Inspected : 1/10/2023"
$ws.Range("C18").Value = "Date Received : 1/10/2023"

# --- Swap the Completed / Partial markers (R <-> *) ---
$ws.Range("C19").Value = "*"
$ws.Range("C20").Value = "R"

# --- Insert a new row for "Partial quantity" and fill it in ---
$ws.Rows(23).Insert()
$ws.Range("D21").Value = "Partial quantity: 0"
